# Scheduled runner update: refresh cached market-price derived columns
# (currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) on the
# Leve-profit sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) for the rows whose
# underlying market data changed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1419.55
$ws.Range("I28").Value = 1324.5625
$ws.Range("J28").Value = 1799.5
$ws.Range("K28").Value = 1324.5625
$ws.Range("L28").Value = 1799.5
$ws.Range("M28").Value = -839.5625
$ws.Range("N28").Value = -2769.5

$ws.Range("H69").Value = 5214.2856
$ws.Range("I69").Value = 4000
$ws.Range("K69").Value = 12000
$ws.Range("M69").Value = -11126

$ws.Range("H72").Value = 5214.2856
$ws.Range("I72").Value = 4000
$ws.Range("K72").Value = 36000
$ws.Range("M72").Value = -31632

$ws.Range("H86").Value = 2334.6365
$ws.Range("I86").Value = 1909
$ws.Range("J86").Value = 4250
$ws.Range("K86").Value = 1909
$ws.Range("L86").Value = 4250
$ws.Range("M86").Value = -786
$ws.Range("N86").Value = -6496

$ws.Range("H89").Value = 2334.6365
$ws.Range("I89").Value = 1909
$ws.Range("J89").Value = 4250
$ws.Range("K89").Value = 9545
$ws.Range("L89").Value = 21250
$ws.Range("M89").Value = -3929
$ws.Range("N89").Value = -32482

$ws.Range("H100").Value = 2899.2144
$ws.Range("I100").Value = 2208.9
$ws.Range("K100").Value = 2208.9
$ws.Range("M100").Value = -1667.9

$ws.Range("H107").Value = 1811.9412
$ws.Range("I107").Value = 387.1
$ws.Range("J107").Value = 3847.4285
$ws.Range("K107").Value = 387.1
$ws.Range("L107").Value = 3847.4285
$ws.Range("M107").Value = 1532.9
$ws.Range("N107").Value = -7687.4285

$ws.Range("H131").Value = 1117.6786
$ws.Range("I131").Value = 1118.1482
$ws.Range("J131").Value = 1105
$ws.Range("K131").Value = 3354.4446
$ws.Range("L131").Value = 3315
$ws.Range("M131").Value = 1685.5554
$ws.Range("N131").Value = -13395

$ws.Range("H138").Value = 2808.5425
$ws.Range("I138").Value = 1609.0555
$ws.Range("J138").Value = 3335.1462
$ws.Range("K138").Value = 4827.166499999999
$ws.Range("L138").Value = 10005.4386
$ws.Range("M138").Value = 312.8335000000006
$ws.Range("N138").Value = -20285.4386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 1475
$ws.Range("J27").Value = 1475
$ws.Range("L27").Value = 1475
$ws.Range("N27").Value = -1843

$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2616

$ws.Range("H30").Value = 2155.8
$ws.Range("I30").Value = 144.5
$ws.Range("J30").Value = 3496.6667
$ws.Range("K30").Value = 144.5
$ws.Range("L30").Value = 3496.6667
$ws.Range("M30").Value = 5.5
$ws.Range("N30").Value = -3796.6667

$ws.Range("H32").Value = 2921276.8
$ws.Range("I32").Value = 472401
$ws.Range("J32").Value = 19309906
$ws.Range("K32").Value = 472401
$ws.Range("L32").Value = 19309906
$ws.Range("M32").Value = -472114
$ws.Range("N32").Value = -19310480

$ws.Range("H35").Value = 1615.5
$ws.Range("I35").Value = 1615.5
$ws.Range("K35").Value = 1615.5
$ws.Range("M35").Value = -1209.5

$ws.Range("H74").Value = 2286.8667
$ws.Range("I74").Value = 1441.8695
$ws.Range("J74").Value = 5063.2856
$ws.Range("K74").Value = 1441.8695
$ws.Range("L74").Value = 5063.2856
$ws.Range("M74").Value = -567.8695
$ws.Range("N74").Value = -6811.2856

$ws.Range("H77").Value = 2286.8667
$ws.Range("I77").Value = 1441.8695
$ws.Range("J77").Value = 5063.2856
$ws.Range("K77").Value = 7209.3475
$ws.Range("L77").Value = 25316.428
$ws.Range("M77").Value = -2841.3475
$ws.Range("N77").Value = -34052.428

$ws.Range("H122").Value = 2277
$ws.Range("I122").Value = 2277
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6831
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4381
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4883
$ws.Range("I132").Value = 4876.1177
$ws.Range("K132").Value = 14628.3531
$ws.Range("M132").Value = -12098.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1166.6666
$ws.Range("I22").Value = 750
$ws.Range("K22").Value = 750
$ws.Range("M22").Value = -400

$ws.Range("H31").Value = 2846.9167
$ws.Range("I31").Value = 1528.44
$ws.Range("K31").Value = 1528.44
$ws.Range("M31").Value = -1233.44

$ws.Range("H34").Value = 2846.9167
$ws.Range("I34").Value = 1528.44
$ws.Range("K34").Value = 1528.44
$ws.Range("M34").Value = -1326.44

$ws.Range("H58").Value = 979.6
$ws.Range("I58").Value = 974.25
$ws.Range("K58").Value = 974.25
$ws.Range("M58").Value = -771.25

$ws.Range("H70").Value = 22222
$ws.Range("J70").Value = 22222
$ws.Range("L70").Value = 22222
$ws.Range("N70").Value = -22852

$ws.Range("H73").Value = 22222
$ws.Range("J73").Value = 22222
$ws.Range("L73").Value = 22222
$ws.Range("N73").Value = -24406

$ws.Range("H134").Value = 906
$ws.Range("I134").Value = 906
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2718
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -183
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 979.6
$ws.Range("I136").Value = 974.25
$ws.Range("K136").Value = 2922.75
$ws.Range("M136").Value = -372.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5406
$ws.Range("J70").Value = 5971.6665
$ws.Range("L70").Value = 17914.9995
$ws.Range("N70").Value = -18544.9995

$ws.Range("H73").Value = 5406
$ws.Range("J73").Value = 5971.6665
$ws.Range("L73").Value = 17914.9995
$ws.Range("N73").Value = -20098.9995

$ws.Range("H107").Value = 769.44446
$ws.Range("I107").Value = 736
$ws.Range("J107").Value = 779
$ws.Range("K107").Value = 2208
$ws.Range("L107").Value = 2337
$ws.Range("M107").Value = -288
$ws.Range("N107").Value = -6177

$ws.Range("H131").Value = 19321872
$ws.Range("I131").Value = 41751004
$ws.Range("K131").Value = 125253012
$ws.Range("M131").Value = -125247972

$ws.Range("H132").Value = 1612.2858
$ws.Range("J132").Value = 1686.0714
$ws.Range("L132").Value = 15174.6426
$ws.Range("N132").Value = -20234.6426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9778.532999999999
$ws.Range("I80").Value = 13415.889
$ws.Range("K80").Value = 13415.889
$ws.Range("M80").Value = -12417.889

$ws.Range("H83").Value = 9778.532999999999
$ws.Range("I83").Value = 13415.889
$ws.Range("K83").Value = 67079.44499999999
$ws.Range("M83").Value = -62087.44499999999

$ws.Range("H122").Value = 1297.5
$ws.Range("I122").Value = 1297.5
$ws.Range("K122").Value = 3892.5
$ws.Range("M122").Value = -1442.5

$ws.Range("H123").Value = 54757.082
$ws.Range("J123").Value = 54757.082
$ws.Range("L123").Value = 54757.082
$ws.Range("N123").Value = -59657.082

$ws.Range("H132").Value = 4443.0356
$ws.Range("I132").Value = 4433.8335
$ws.Range("K132").Value = 13301.5005
$ws.Range("M132").Value = -10771.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5323.6665
$ws.Range("I82").Value = 5500
$ws.Range("J82").Value = 4971
$ws.Range("K82").Value = 5500
$ws.Range("L82").Value = 4971
$ws.Range("M82").Value = -5139
$ws.Range("N82").Value = -5693

$ws.Range("H85").Value = 5323.6665
$ws.Range("I85").Value = 5500
$ws.Range("J85").Value = 4971
$ws.Range("K85").Value = 5500
$ws.Range("L85").Value = 4971
$ws.Range("M85").Value = -4252
$ws.Range("N85").Value = -7467

$ws.Range("H128").Value = 29995
$ws.Range("J128").Value = 29995
$ws.Range("L128").Value = 29995
$ws.Range("N128").Value = -39955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 30346.2
$ws.Range("I45").Value = 19998
$ws.Range("J45").Value = 31496
$ws.Range("K45").Value = 19998
$ws.Range("L45").Value = 31496
$ws.Range("M45").Value = -19507
$ws.Range("N45").Value = -32478

$ws.Range("H100").Value = 5683.1665
$ws.Range("I100").Value = 6556.6
$ws.Range("J100").Value = 1316
$ws.Range("K100").Value = 13113.2
$ws.Range("L100").Value = 2632
$ws.Range("M100").Value = -12572.2
$ws.Range("N100").Value = -3714
